$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41, shifting existing rows 41:56 down to 42:57.
$ws.Rows.Item(41).Insert()

# Populate the newly-inserted row 41 with the new weekly record.
$ws.Range("A41").Value = 5
$ws.Range("B41").Value = "Macroferia Regional de Talca"
$ws.Range("C41").Value = "Maule"
$ws.Range("D41").Value = 44508
$ws.Range("E41").Value = 7
$ws.Range("F41").Value = 100112022
$ws.Range("G41").Value = "Arveja Verde"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 600
$ws.Range("K41").Value = 12000
$ws.Range("L41").Value = 12000
$ws.Range("M41").Value = 12000
$ws.Range("N41").Value = "$/saco 25 kilos"
$ws.Range("O41").Value = "Región del Maule"
$ws.Range("P41").Value = 480
$ws.Range("Q41").Value = 25
$ws.Range("R41").Value = "Hortaliza"
